$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Proximity")

$data = @(
    @("2026-02-01", "13:16:44", "13:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "13:16:44", "13:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "13:16:46", "13:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "13:17:11", "13:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "13:17:43", "13:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door")
)

$startRow = 5
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 1
        $cell = $ws.Cells.Item($row, $col)
        # Column A holds date-looking text ("2026-02-01"); force text entry so
        # Excel doesn't silently convert it to a date serial number, then
        # restore the default "Normal" style so no stray number format sticks.
        if ($col -eq 1) {
            $cell.NumberFormat = "@"
            $cell.Value = $rowData[$c]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $rowData[$c]
        }
    }
}
